# Apply updated "想去人数" (want-to-go count) values in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet updates (column F, row -> new value) ---
$wsExhibition.Cells.Item(2, 6).Value = 55
$wsExhibition.Cells.Item(3, 6).Value = 3216
$wsExhibition.Cells.Item(4, 6).Value = 1683
$wsExhibition.Cells.Item(5, 6).Value = 2259
$wsExhibition.Cells.Item(7, 6).Value = 319
$wsExhibition.Cells.Item(8, 6).Value = 1188
$wsExhibition.Cells.Item(9, 6).Value = 1034
$wsExhibition.Cells.Item(10, 6).Value = 256
$wsExhibition.Cells.Item(11, 6).Value = 473
$wsExhibition.Cells.Item(16, 6).Value = 7990
$wsExhibition.Cells.Item(17, 6).Value = 349
$wsExhibition.Cells.Item(19, 6).Value = 224
$wsExhibition.Cells.Item(20, 6).Value = 237
$wsExhibition.Cells.Item(21, 6).Value = 170
$wsExhibition.Cells.Item(23, 6).Value = 545
$wsExhibition.Cells.Item(25, 6).Value = 1137
$wsExhibition.Cells.Item(26, 6).Value = 986
$wsExhibition.Cells.Item(27, 6).Value = 1701
$wsExhibition.Cells.Item(28, 6).Value = 205
$wsExhibition.Cells.Item(29, 6).Value = 62
$wsExhibition.Cells.Item(30, 6).Value = 1680
$wsExhibition.Cells.Item(31, 6).Value = 239
$wsExhibition.Cells.Item(33, 6).Value = 479
$wsExhibition.Cells.Item(39, 6).Value = 190
$wsExhibition.Cells.Item(40, 6).Value = 355
$wsExhibition.Cells.Item(42, 6).Value = 226

# --- 全部类型 sheet updates (column F, row -> new value) ---
$wsAllTypes.Cells.Item(4, 6).Value = 55
$wsAllTypes.Cells.Item(5, 6).Value = 3216
$wsAllTypes.Cells.Item(6, 6).Value = 1683
$wsAllTypes.Cells.Item(7, 6).Value = 2259
$wsAllTypes.Cells.Item(9, 6).Value = 319
$wsAllTypes.Cells.Item(10, 6).Value = 1188
$wsAllTypes.Cells.Item(12, 6).Value = 1034
$wsAllTypes.Cells.Item(13, 6).Value = 256
$wsAllTypes.Cells.Item(14, 6).Value = 473
$wsAllTypes.Cells.Item(18, 6).Value = 7991
$wsAllTypes.Cells.Item(19, 6).Value = 349
$wsAllTypes.Cells.Item(22, 6).Value = 224
$wsAllTypes.Cells.Item(23, 6).Value = 237
$wsAllTypes.Cells.Item(24, 6).Value = 170
$wsAllTypes.Cells.Item(26, 6).Value = 545
$wsAllTypes.Cells.Item(28, 6).Value = 1137
$wsAllTypes.Cells.Item(29, 6).Value = 986
$wsAllTypes.Cells.Item(30, 6).Value = 1703
$wsAllTypes.Cells.Item(31, 6).Value = 205
$wsAllTypes.Cells.Item(32, 6).Value = 62
$wsAllTypes.Cells.Item(33, 6).Value = 1680
$wsAllTypes.Cells.Item(34, 6).Value = 239
$wsAllTypes.Cells.Item(36, 6).Value = 479
$wsAllTypes.Cells.Item(42, 6).Value = 190
$wsAllTypes.Cells.Item(43, 6).Value = 355
$wsAllTypes.Cells.Item(49, 6).Value = 226

